$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "ODI Batting": drop the stray empty INNING_NUMBER cells on rows 4/6/8/10
#    (they were placeholders scraped with no value - the edit removes them
#    entirely rather than leaving an empty cell behind).
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("B4").ClearContents()
$battingSheet.Range("B6").ClearContents()
$battingSheet.Range("B8").ClearContents()
$battingSheet.Range("B10").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Batting Extra" sheet (extra scraped batting fields),
#    placed after the existing sheets.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Match the outline / page-margin conventions already used by the sibling
# sheets in this workbook.
$extra.Outline.SummaryRow = 1
$extra.Outline.SummaryColumn = 1
$extra.PageSetup.LeftMargin = 54
$extra.PageSetup.RightMargin = 54
$extra.PageSetup.TopMargin = 72
$extra.PageSetup.BottomMargin = 72
$extra.PageSetup.HeaderMargin = 36
$extra.PageSetup.FooterMargin = 36

# Header row - reuse the same bold/centered/bordered header style already
# used by the other sheets in this workbook (copy its format onto the new
# header cells rather than building a brand-new style).
$headerStyleSource = $battingSheet.Range("A1")
$headerStyleSource.Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $extra.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# Numeric-looking text (MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL) is
# entered with a leading apostrophe so Excel keeps it as text instead of
# coercing it to a number; BATTING_POSITION is a real number; blank entries
# are left as "" (scrape gaps) to match the source cells that were empty.
$rows = @(
    @("4100", 7, "2", "0", "6.02%",  "NO"),
    @("4101", 7, "1", "2", "15.08%", "NO"),
    @("4102", "", "", "", "", "NO"),
    @("4103", 9, "0", "0", "", "NO"),
    @("4104", 8, "", "", "", "NO"),
    @("4105", 8, "1", "0", "1.95%", "NO"),
    @("4248", 9, "", "", "", "NO"),
    @("4249", 9, "1", "0", "4.61%", "NO"),
    @("4251", 9, "", "", "", "NO")
)

$r = 2
foreach ($row in $rows) {
    $extra.Cells.Item($r, 1).Value = "'" + $row[0]

    if ($row[1] -eq "") {
        $extra.Cells.Item($r, 2).Value = "'"
    } else {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }

    for ($col = 3; $col -le 5; $col++) {
        $val = $row[$col - 1]
        if ($val -eq "") {
            $extra.Cells.Item($r, $col).Value = "'"
        } else {
            $extra.Cells.Item($r, $col).Value = "'" + $val
        }
    }

    $extra.Cells.Item($r, 6).Value = $row[5]

    $r++
}

# Leave the original first sheet selected/active, same as before the edit.
$wb.Worksheets.Item(1).Activate()
